# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date", "Correspond Handoff Datetime" and
# "Correspond Handback DateTime" timestamps for the 7827cd75-... entry across
# the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-09-09 12:24:44"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-09-09 12:24:33"
$zhcn.Range("K3").Value = "2016-09-09 12:25:49"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-09-09 12:24:44"
$dede.Range("K3").Value = "2016-09-09 12:26:15"
